$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Recolor the existing "Invoervelden Uitgave" / "Overzicht uitgaven"
#    block with the green accent font (matches RGB 146,208,80 = 0x92D050)
# ---------------------------------------------------------------------
$green = 146 + (208 * 256) + (80 * 65536)

$ws.Range("A3:B6").Font.Color  = $green
$ws.Range("A8").Font.Color     = $green
$ws.Range("A11:M14").Font.Color = $green

# ---------------------------------------------------------------------
# 2) New "ToDo" list starting at row 16
# ---------------------------------------------------------------------
$ws.Range("A16").Value = "ToDo"

$ws.Range("A17").Value = "basics:"
$ws.Range("A17").Font.Bold = $true

$ws.Range("A18").Value = "input form maken"
$ws.Range("G18").Value = "Daan"

$ws.Range("A19").Value = 'veld toevoegen en "ingevoerd door"'
$ws.Range("G19").Value = "Daan"

$ws.Range("A20").Value = "total amount pp maken"

$ws.Range("A21").Value = "delete buttons toevoegen"
$ws.Range("G21").Value = "Tim"

$ws.Range("A22").Value = "sort buttons wijzigen naar dropdown"
$ws.Range("G22").Value = "Pair"

$ws.Range("A24").Value = "werkend maken:"
$ws.Range("A24").Font.Bold = $true

$ws.Range("A25").Value = "input form werkend maken"
$ws.Range("A25").Font.Bold = $false
$ws.Range("G25").Value = "Daan"

$ws.Range("A26").Value = "sort buttons werkend maken"
$ws.Range("A26").Font.Bold = $false
$ws.Range("G26").Value = "Pair"

$ws.Range("A27").Value = "total amount pp werkend maken"
$ws.Range("A27").Font.Bold = $false

$ws.Range("A28").Value = "delete buttons werkend maken"
$ws.Range("A28").Font.Bold = $false
$ws.Range("G28").Value = "Tim"

$ws.Range("A30").Value = "extra:"
$ws.Range("A30").Font.Bold = $true

$ws.Range("A31").Value = "communicatie met API ipv data.js"
$ws.Range("A31").Font.Bold = $false

$ws.Range("A32").Value = "user accounts maken met inlogsysteem"
$ws.Range("A32").Font.Bold = $false

$ws.Range("A33").Value = "styling"
$ws.Range("A33").Font.Bold = $false

$ws.Range("A35").Value = "final:"
$ws.Range("A35").Font.Bold = $true

$ws.Range("A36").Value = "hosting"
$ws.Range("A36").Font.Bold = $false

# ---------------------------------------------------------------------
# 3) Selection as left by the author
# ---------------------------------------------------------------------
$ws.Range("G26").Select()
